$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Angptl1"
$ws.Range("C2").Value = "Tek"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 55.82849233333334
$ws.Range("H2").Value = 167.485477
$ws.Range("I2").Value = 0.9717741676025277
$ws.Range("J2").Value = 0.9717741676025275
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 44.69746666666666
$ws.Range("N2").Value = 134.0924
$ws.Range("O2").Value = 0.6823972194925493
$ws.Range("P2").Value = 0.6823972194925493
$ws.Range("Q2").Value = 2495.392175119422
$ws.Range("R2").Value = 22458.5295760748
$ws.Range("S2").Value = 0.6631359899466515
$ws.Range("T2").Value = 0.6631359899466513

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Angptl1"
$ws.Range("C3").Value = "Tek"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 55.82849233333334
$ws.Range("H3").Value = 167.485477
$ws.Range("I3").Value = 0.9717741676025277
$ws.Range("J3").Value = 0.9717741676025275
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.63243533333333
$ws.Range("N3").Value = 55.897306
$ws.Range("O3").Value = 0.2844618053784121
$ws.Range("P3").Value = 0.2844618053784121
$ws.Range("Q3").Value = 1040.220773158329
$ws.Range("R3").Value = 9361.986958424963
$ws.Range("S3").Value = 0.2764326341363187
$ws.Range("T3").Value = 0.2764326341363186

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Angptl1"
$ws.Range("C4").Value = "Tek"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 55.82849233333334
$ws.Range("H4").Value = 167.485477
$ws.Range("I4").Value = 0.9717741676025277
$ws.Range("J4").Value = 0.9717741676025275
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.170755666666667
$ws.Range("N4").Value = 6.512267
$ws.Range("O4").Value = 0.03314097512903853
$ws.Range("P4").Value = 0.03314097512903853
$ws.Range("Q4").Value = 121.1900160940399
$ws.Range("R4").Value = 1090.710144846359
$ws.Range("S4").Value = 0.03220554351955749
$ws.Range("T4").Value = 0.03220554351955748

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Angptl1"
$ws.Range("C5").Value = "Tek"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.621576
$ws.Range("H5").Value = 4.864728
$ws.Range("I5").Value = 0.02822583239747235
$ws.Range("J5").Value = 0.02822583239747234
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 44.69746666666666
$ws.Range("N5").Value = 134.0924
$ws.Range("O5").Value = 0.6823972194925493
$ws.Range("P5").Value = 0.6823972194925493
$ws.Range("Q5").Value = 72.48033920746667
$ws.Range("R5").Value = 652.3230528672001
$ws.Range("S5").Value = 0.01926122954589785
$ws.Range("T5").Value = 0.01926122954589784

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Angptl1"
$ws.Range("C6").Value = "Tek"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.621576
$ws.Range("H6").Value = 4.864728
$ws.Range("I6").Value = 0.02822583239747235
$ws.Range("J6").Value = 0.02822583239747234
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.63243533333333
$ws.Range("N6").Value = 55.897306
$ws.Range("O6").Value = 0.2844618053784121
$ws.Range("P6").Value = 0.2844618053784121
$ws.Range("Q6").Value = 30.21390995808533
$ws.Range("R6").Value = 271.925189622768
$ws.Range("S6").Value = 0.008029171242093458
$ws.Range("T6").Value = 0.008029171242093458

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Angptl1"
$ws.Range("C7").Value = "Tek"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.621576
$ws.Range("H7").Value = 4.864728
$ws.Range("I7").Value = 0.02822583239747235
$ws.Range("J7").Value = 0.02822583239747234
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.170755666666667
$ws.Range("N7").Value = 6.512267
$ws.Range("O7").Value = 0.03314097512903853
$ws.Range("P7").Value = 0.03314097512903853
$ws.Range("Q7").Value = 3.520045290930668
$ws.Range("R7").Value = 31.68040761837601
$ws.Range("S7").Value = 0.0009354316094810411
$ws.Range("T7").Value = 0.0009354316094810407
